$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply style to new header cell W1 (copy format from V1, then set value)
$ws.Range("V1").Copy($ws.Range("W1")) | Out-Null
$ws.Range("W1").Value = 21

# Apply style to new index cells A8:A11 (copy format from A7), then set values
$ws.Range("A7").Copy($ws.Range("A8")) | Out-Null
$ws.Range("A7").Copy($ws.Range("A9")) | Out-Null
$ws.Range("A7").Copy($ws.Range("A10")) | Out-Null
$ws.Range("A7").Copy($ws.Range("A11")) | Out-Null
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Row 2: HKL header labels (set all to ensure correctness)
$ws.Range("C2").Value = '[1, 1, 1]'
$ws.Range("D2").Value = '[2, 0, 0]'
$ws.Range("E2").Value = '[2, 2, 0]'
$ws.Range("F2").Value = '[3, 1, 1]'
$ws.Range("G2").Value = '[2, 2, 2]'
$ws.Range("H2").Value = '[4, 0, 0]'
$ws.Range("I2").Value = '[3, 3, 1]'
$ws.Range("J2").Value = '[4, 2, 0]'
$ws.Range("K2").Value = '[4, 2, 2]'
$ws.Range("L2").Value = '[5, 1, 1]'
$ws.Range("M2").Value = '[3, 3, 3]'
$ws.Range("N2").Value = '1Pair-A'
$ws.Range("O2").Value = '1Pair-B'
$ws.Range("P2").Value = '2Pairs-A'
$ws.Range("Q2").Value = '2Pairs-B'
$ws.Range("R2").Value = '3Pairs-A'
$ws.Range("S2").Value = '3Pairs-B'
$ws.Range("T2").Value = '3Pairs-C'
$ws.Range("U2").Value = '4Pairs'
$ws.Range("V2").Value = '5A4F'
$ws.Range("W2").Value = 'MaxUnique'

# Column B labels for rows 3-11
$ws.Range("B3").Value = 'Equal Angle'
$ws.Range("B4").Value = 'CLR'
$ws.Range("B5").Value = 'BT8Hex'
$ws.Range("B6").Value = 'Spiral'
$ws.Range("B7").Value = 'OffsetF'
$ws.Range("B8").Value = 'OffsetA'
$ws.Range("B9").Value = 'RD Single'
$ws.Range("B10").Value = 'TD Single'
$ws.Range("B11").Value = 'HexGrid-90degTilt5degRes'

# Numeric data cells C:W for rows 3-11
$ws.Range("C3").Value = 1.050158501440922
$ws.Range("D3").Value = 0.8878458213256484
$ws.Range("E3").Value = 1.02193083573487
$ws.Range("F3").Value = 0.9745749279538904
$ws.Range("G3").Value = 1.050158501440922
$ws.Range("H3").Value = 0.8878458213256484
$ws.Range("I3").Value = 1.029517291066282
$ws.Range("J3").Value = 0.976707492795389
$ws.Range("K3").Value = 1.019279538904899
$ws.Range("L3").Value = 0.9312608069164265
$ws.Range("M3").Value = 1.050151296829971
$ws.Range("N3").Value = 1.050158501440922
$ws.Range("O3").Value = 1.02193083573487
$ws.Range("P3").Value = 0.9548883285302594
$ws.Range("Q3").Value = 0.9982528818443804
$ws.Range("R3").Value = 0.9866450528338137
$ws.Range("S3").Value = 0.9614505283381364
$ws.Range("T3").Value = 0.9866450528338137
$ws.Range("U3").Value = 0.9836275216138328
$ws.Range("V3").Value = 0.9969337175792508
$ws.Range("W3").Value = 0.9864094020172911

$ws.Range("C4").Value = 1.00369818621128
$ws.Range("D4").Value = 0.984758933856594
$ws.Range("E4").Value = 0.9939828252518585
$ws.Range("F4").Value = 0.9915119850352294
$ws.Range("G4").Value = 1.00369818621128
$ws.Range("H4").Value = 0.984758933856594
$ws.Range("I4").Value = 0.996148678369443
$ws.Range("J4").Value = 0.9922747580120095
$ws.Range("K4").Value = 0.9973039466221253
$ws.Range("L4").Value = 0.9866765953859452
$ws.Range("M4").Value = 1.003685926047098
$ws.Range("N4").Value = 1.00369818621128
$ws.Range("O4").Value = 0.9939828252518585
$ws.Range("P4").Value = 0.9893708795542262
$ws.Range("Q4").Value = 0.992747405143544
$ws.Range("R4").Value = 0.9941466484399109
$ws.Range("S4").Value = 0.9900845813812272
$ws.Range("T4").Value = 0.994146648439911
$ws.Range("U4").Value = 0.9934879825887406
$ws.Range("V4").Value = 0.9955300233132485
$ws.Range("W4").Value = 0.9932944885930607

$ws.Range("C5").Value = 1.007341267463998
$ws.Range("D5").Value = 0.9769087892911116
$ws.Range("E5").Value = 0.9956190279053817
$ws.Range("F5").Value = 0.9904191524181927
$ws.Range("G5").Value = 1.007341267463998
$ws.Range("H5").Value = 0.9769087892911116
$ws.Range("I5").Value = 0.9981166925210559
$ws.Range("J5").Value = 0.9914007440302651
$ws.Range("K5").Value = 0.9989154577641843
$ws.Range("L5").Value = 0.9823576474895908
$ws.Range("M5").Value = 1.007336105493718
$ws.Range("N5").Value = 1.007341267463998
$ws.Range("O5").Value = 0.9956190279053817
$ws.Range("P5").Value = 0.9862639085982466
$ws.Range("Q5").Value = 0.9930190901617872
$ws.Range("R5").Value = 0.9932896948868303
$ws.Range("S5").Value = 0.987648989871562
$ws.Range("T5").Value = 0.9932896948868303
$ws.Range("U5").Value = 0.992572059269671
$ws.Range("V5").Value = 0.9955259009085363
$ws.Range("W5").Value = 0.9926348473604725

$ws.Range("C6").Value = 0.9964632770323403
$ws.Range("D6").Value = 0.9924118810084416
$ws.Range("E6").Value = 0.993773031237944
$ws.Range("F6").Value = 0.9930489074098244
$ws.Range("G6").Value = 0.9964632770323403
$ws.Range("H6").Value = 0.9924118810084416
$ws.Range("I6").Value = 0.99411028412095
$ws.Range("J6").Value = 0.9940568388390248
$ws.Range("K6").Value = 0.9949920963782543
$ws.Range("L6").Value = 0.9915570012361647
$ws.Range("M6").Value = 0.996449271428874
$ws.Range("N6").Value = 0.9964632770323403
$ws.Range("O6").Value = 0.993773031237944
$ws.Range("P6").Value = 0.9930924561231929
$ws.Range("Q6").Value = 0.9934109693238842
$ws.Range("R6").Value = 0.9942160630929088
$ws.Range("S6").Value = 0.9930779398854034
$ws.Range("T6").Value = 0.9942160630929088
$ws.Range("U6").Value = 0.9939242741721377
$ws.Range("V6").Value = 0.9944320747441783
$ws.Range("W6").Value = 0.993801664657868

$ws.Range("C7").Value = 1.102512415797626
$ws.Range("D7").Value = 1.100149909827876
$ws.Range("E7").Value = 0.8635026292527483
$ws.Range("F7").Value = 1.009573608983402
$ws.Range("G7").Value = 1.102512415797626
$ws.Range("H7").Value = 1.100149909827876
$ws.Range("I7").Value = 0.945893760916121
$ws.Range("J7").Value = 0.9475770775616994
$ws.Range("K7").Value = 1.024563217950286
$ws.Range("L7").Value = 1.044359231486226
$ws.Range("M7").Value = 1.102512415797626
$ws.Range("N7").Value = 1.102512415797626
$ws.Range("O7").Value = 0.8635026292527483
$ws.Range("P7").Value = 0.9818262695403119
$ws.Range("Q7").Value = 0.9365381191180749
$ws.Range("R7").Value = 1.022054984959417
$ws.Range("S7").Value = 0.9910753826880084
$ws.Range("T7").Value = 1.022054984959417
$ws.Range("U7").Value = 1.018934640965413
$ws.Range("V7").Value = 1.035650195931856
$ws.Range("W7").Value = 1.004766481471998

$ws.Range("C8").Value = 0.9351743044627683
$ws.Range("D8").Value = 0.9625785200588551
$ws.Range("E8").Value = 1.045963227266105
$ws.Range("F8").Value = 0.9924926841183469
$ws.Range("G8").Value = 0.9351743044627683
$ws.Range("H8").Value = 0.9625785200588551
$ws.Range("I8").Value = 1.005690499239094
$ws.Range("J8").Value = 1.014279601002294
$ws.Range("K8").Value = 0.9784386172537788
$ws.Range("L8").Value = 0.9868743649618772
$ws.Range("M8").Value = 0.9351743044627683
$ws.Range("N8").Value = 0.9351743044627683
$ws.Range("O8").Value = 1.045963227266105
$ws.Range("P8").Value = 1.00427087366248
$ws.Range("Q8").Value = 1.019227955692226
$ws.Range("R8").Value = 0.9812386839292428
$ws.Range("S8").Value = 1.000344810481102
$ws.Range("T8").Value = 0.9812386839292428
$ws.Range("U8").Value = 0.9840521839765188
$ws.Range("V8").Value = 0.9742766080737688
$ws.Range("W8").Value = 0.9901864772953899

$ws.Range("C9").Value = 1.88
$ws.Range("D9").Value = 0.23
$ws.Range("E9").Value = 0.89
$ws.Range("F9").Value = 0.85
$ws.Range("G9").Value = 1.88
$ws.Range("H9").Value = 0.23
$ws.Range("I9").Value = 1.2
$ws.Range("J9").Value = 0.6899999999999999
$ws.Range("K9").Value = 1.33
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 1.88
$ws.Range("N9").Value = 1.88
$ws.Range("O9").Value = 0.89
$ws.Range("P9").Value = 0.5600000000000001
$ws.Range("Q9").Value = 0.87
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 0.6566666666666667
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 0.9625
$ws.Range("V9").Value = 1.146
$ws.Range("W9").Value = 0.9462499999999999

$ws.Range("C10").Value = 1.63
$ws.Range("D10").Value = 0.22
$ws.Range("E10").Value = 1.14
$ws.Range("F10").Value = 0.8100000000000001
$ws.Range("G10").Value = 1.63
$ws.Range("H10").Value = 0.22
$ws.Range("I10").Value = 1.29
$ws.Range("J10").Value = 0.78
$ws.Range("K10").Value = 1.23
$ws.Range("L10").Value = 0.47
$ws.Range("M10").Value = 1.63
$ws.Range("N10").Value = 1.63
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 0.6799999999999999
$ws.Range("Q10").Value = 0.975
$ws.Range("R10").Value = 0.9966666666666666
$ws.Range("S10").Value = 0.7233333333333333
$ws.Range("T10").Value = 0.9966666666666666
$ws.Range("U10").Value = 0.95
$ws.Range("V10").Value = 1.086
$ws.Range("W10").Value = 0.94625

$ws.Range("C11").Value = 0.9943153340810815
$ws.Range("D11").Value = 0.9939049416825618
$ws.Range("E11").Value = 0.9941624606196982
$ws.Range("F11").Value = 0.9934084552270573
$ws.Range("G11").Value = 0.9943153340810815
$ws.Range("H11").Value = 0.9939049416825618
$ws.Range("I11").Value = 0.9935580259769844
$ws.Range("J11").Value = 0.9950119991389479
$ws.Range("K11").Value = 0.9942616629124984
$ws.Range("L11").Value = 0.9925809888011333
$ws.Range("M11").Value = 0.9943020649766697
$ws.Range("N11").Value = 0.9943153340810815
$ws.Range("O11").Value = 0.9941624606196982
$ws.Range("P11").Value = 0.99403370115113
$ws.Range("Q11").Value = 0.9937854579233778
$ws.Range("R11").Value = 0.9941275787944471
$ws.Range("S11").Value = 0.9938252858431058
$ws.Range("T11").Value = 0.9941275787944471
$ws.Range("U11").Value = 0.9939477979025997
$ws.Range("V11").Value = 0.994021305138296
$ws.Range("W11").Value = 0.9939004835549954
